$wb = $excel.ActiveWorkbook

# --- Metrics sheet: update the raw metric values in column B ---
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Activate()

$metrics.Range("B2").Value = 92105.53
$metrics.Range("B3").Value = 78431.930000000008
$metrics.Range("B4").Value = 36524.410000000003
$metrics.Range("B5").Value = 3563
$metrics.Range("B6").Value = 4459237.0000000009
$metrics.Range("B7").Value = 3768250.5999999996
$metrics.Range("B8").Value = 1307126.55
$metrics.Range("B9").Value = 172564
$metrics.Range("B10").Value = 32924560.800999828
$metrics.Range("B11").Value = 31043472.120000005
$metrics.Range("B12").Value = 11588835.439999999
$metrics.Range("B13").Value = 1270191

# Move the Metrics sheet's selection to D9
$metrics.Range("D9").Select()

# --- today sheet: move the active selection to G7 ---
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("G7").Select()
